# Add a new "DefaultHitTime" row (row 31) to the Property sheet, matching
# the pattern used by the existing rows (e.g. row 29 "AtkDis").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text columns: Id, Type, Friend (relation), Desc -> force text format like
# the other rows in this table (numFmtId 49 / "@"), then set the value.
$ws.Range("A31").NumberFormat = "@"
$ws.Range("A31").Value = "DefaultHitTime"

$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "float"

$ws.Range("C31").Value = $false
$ws.Range("D31").Value = $false
$ws.Range("E31").Value = $false
$ws.Range("F31").Value = $true

$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0

$ws.Range("I31").NumberFormat = "@"
$ws.Range("I31").Value = "Friend"

$ws.Range("J31").NumberFormat = "@"
$ws.Range("J31").Value = "缺省打击时间（本来应该打到但是物理没碰撞到或者其他原因）"

# The existing "TRUE,FALSE" list validation on column F already spans
# F2:F1048576, so it automatically covers the newly added F31 cell too.
